# Actualización control 3 por nueva versión de Quarto.
# The "Slides class_10" material link had been placed in the
# "material_futuro" (future material) column (I12); move it into the
# current "material" column (G12) now that the class has occurred.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$materialFuturo = $ws.Range("I12")
$material = $ws.Range("G12")

$material.Value = $materialFuturo.Value2
$materialFuturo.ClearContents()

# Update the saved cell selection/active cell on the sheet.
$ws.Range("H17").Select()
